# Adds the "Negative" test-case rows (row 3 and row 4) to each of the three
# sheets (ServicesCategory, ServiceType, ServiceNameEBP), mirroring the
# existing row 2 pattern but for the new "testT4116_Negative" / "NG1" case.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: ServicesCategory
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A3").Value = "testT4116_Negative"
$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = "click"
$ws1.Range("E3").Value = "autoText"
$ws1.Range("F3").Value = "click"
$ws1.Range("G3").Value = "autoText"
$ws1.Range("H3").Value = "autoText"

$ws1.Range("A4").Value = "testT4116_Negative"
$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = "NG1"
$ws1.Range("D4").Value = "click"
$ws1.Range("F4").Value = "click"
$ws1.Range("G4").Value = "autoText"

$ws1.Rows.Item(3).RowHeight = 15
$ws1.Rows.Item(4).RowHeight = 15

$ws1.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------
# Sheet 2: ServiceType
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = "testT4116_Negative"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = "click"
$ws2.Range("E3").Value = "click"
$ws2.Range("F3").Value = "autoText"
$ws2.Range("G3").Value = "click"
$ws2.Range("H3").Value = "autoText"
$ws2.Range("I3").Value = "autoText"

$ws2.Range("A4").Value = "testT4116_Negative"
$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = "NG1"
$ws2.Range("D4").Value = "click"
$ws2.Range("E4").Value = "click"
$ws2.Range("G4").Value = "click"
$ws2.Range("H4").Value = "autoText"

# ---------------------------------------------------------------------
# Sheet 3: ServiceNameEBP
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A3").Value = "testT4116_Negative"
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 1
$ws3.Range("D3").Value = "click"
$ws3.Range("E3").Value = "autoText"
$ws3.Range("F3").Value = "click"
$ws3.Range("G3").Value = "autoText"
$ws3.Range("H3").Value = "No Evidence"
$ws3.Range("I3").Value = "Concerning Practice"
$ws3.Range("J3").Value = "autoText"

$ws3.Range("A4").Value = "testT4116_Negative"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = "'NG1"
$ws3.Range("D4").Value = "click"
$ws3.Range("F4").Value = "click"
$ws3.Range("G4").Value = "autoText"
$ws3.Range("H4").Value = "No Evidence"
$ws3.Range("I4").Value = "Concerning Practice"

$ws3.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: ServiceType and ServiceNameEBP end
# up on C4, ServicesCategory (now the active tab) ends up on E7.
# ---------------------------------------------------------------------
[void]$ws2.Range("C4").Select()
[void]$ws3.Range("C4").Select()

$ws1.Activate()
[void]$ws1.Range("E7").Select()
